# Version 2.0.1 solucionado error espera de base de datos
# Update patient identification data on the "Hoja de Ingreso y Egreso" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nombre del paciente (row 6) ---
$ws.Range("A6").Value = "MACARIO"
$ws.Range("C6").Value = "AGUILAR"
$ws.Range("E6").Value = "OSCAR"
$ws.Range("G6").Value = "ALFREDO"
$ws.Range("I6").Value = "7863/201760947"

# --- Fecha de nacimiento / edad / lugar de nacimiento (row 12) ---
# Values look like a date / a plain number to the parser, so force them to
# stay plain text (matches the original shared-string/text storage) the same
# way a user would by pre-formatting the cell as Text / using a leading
# apostrophe.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "0002-12-08"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "15"
$ws.Range("H12").Value = "GUATEMALA"

# --- Ocupación / nacionalidad / No. de cédula (row 14) ---
$ws.Range("D14").Value = "ESTUDIANTE"
$ws.Range("F14").Value = ""
$ws.Range("H14").Value = "P 472 F 372 L 148"

# --- Contacto de emergencia (row 20) ---
$ws.Range("A20").Value = "ZOILA AGUILAR"
$ws.Range("F20").Value = "MADRE"
$ws.Range("H20").Value = "39CALLE 3AV 39-23 Z.8 "
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "34694016"
